$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell text updates (Coin, Link, Price, Volume(1h) columns)
$updates = [ordered]@{
  'D2' = '61.969.68'
  'E2' = '  -1.40%  '
  'D3' = '3.421.35'
  'E3' = '  -0.32%  '
  'E4' = '  -0.23%  '
  'D5' = '578.95'
  'E5' = '  +0.18%  '
  'D6' = '154.03'
  'E6' = '  +5.25%  '
  'E7' = '  +0.02%  '
  'E8' = '  +1.69%  '
  'E9' = '  +3.29%  '
  'E10' = '  +0.96%  '
  'E11' = '  +3.92%  '
  'D12' = '4.010.70'
  'E12' = '  -0.98%  '
  'E13' = '  +0.75%  '
  'D14' = '28.33'
  'E14' = '  -2.10%  '
  'E15' = '  +0.97%  '
  'D16' = '3.403.96'
  'E16' = '  -1.43%  '
  'D17' = '62.026.70'
  'E17' = '  -1.58%  '
  'D18' = '6.57'
  'E18' = '  +2.95%  '
  'D19' = '14.46'
  'E19' = '  +0.67%  '
  'D20' = '8.95'
  'E20' = '  -2.71%  '
  'D21' = '382.69'
  'E21' = '  -0.33%  '
  'D22' = '0.570'
  'E22' = '  +2.02%  '
  'D23' = '75.94'
  'E23' = '  +2.23%  '
  'E24' = '  +0.07%  '
  'D25' = '3.561.46'
  'E25' = '  -1.10%  '
  'E26' = '  -1.39%  '
  'E27' = '  -0.75%  '
  'D28' = '7.63'
  'E28' = '  +0.47%  '
  'D29' = '0.996'
  'E29' = '  -0.35%  '
  'E30' = '  +0.78%  '
  'E31' = '  -2.75%  '
  'E32' = '  -0.09%  '
  'D33' = '23.30'
  'E33' = '  +0.33%  '
  'E34' = '  +0.24%  '
  'E35' = '  +5.40%  '
  'D36' = '1.60'
  'E36' = '  +0.64%  '
  'D37' = '6.96'
  'E37' = '  -1.83%  '
  'B38' = 'Monero'
  'C38' = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
  'D38' = '168.14'
  'E38' = '  -0.06%  '
  'B39' = 'EnergySwap'
  'C39' = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
  'D39' = '30.94'
  'E39' = '  -2.49%  '
  'D40' = '3.460.00'
  'E40' = '  -0.64%  '
  'E41' = '  +2.67%  '
  'D42' = '42.70'
  'E42' = '  +1.00%  '
  'D43' = '0.780'
  'E43' = '  -0.84%  '
  'E44' = '  +2.20%  '
  'E45' = '  -2.37%  '
  'D46' = '1.17'
  'E46' = '  -2.81%  '
  'D47' = '2.555.95'
  'E47' = '  -1.08%  '
  'B48' = 'InjectiveProtocol'
  'C48' = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
  'D48' = '23.08'
  'E48' = '  +1.40%  '
  'B49' = 'Cosmos'
  'C49' = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
  'D49' = '6.81'
  'E49' = '  +0.61%  '
  'D50' = '2.20'
  'E50' = '  -3.43%  '
  'E51' = '  -0.44%  '
}

# Columns that can contain plain numeric-looking text (e.g. "578.95", "23.30")
# must be forced to Text format first, otherwise Excel silently converts them to
# floating point numbers and mangles exact formatting (trailing zeros, "61.969.68"
# style thousand-grouped strings, etc.)
foreach ($addr in $updates.Keys) {
  if ($addr -match "^D") {
    $ws.Range($addr).NumberFormat = "@"
  }
  $ws.Range($addr).Value = $updates[$addr]
}

Write-Host "Applied $($updates.Count) cell updates"